$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New raw field data for rows 127:132, column H (field/lab measurement)
$ws.Range("H127").Value = 1.3001
$ws.Range("H128").Value = 1.2985
$ws.Range("H129").Value = 1.2917000000000001
$ws.Range("H130").Value = 1.332
$ws.Range("H131").Value = 1.2887
$ws.Range("H132").Value = 1.3153999999999999

# Extend the existing computed-column formulas (J = (F-H)*1000, K = J/(C/1000))
# down through the newly populated rows, matching the pattern already used
# in the column above.
$ws.Range("J127:J132").Formula = "=(F127-H127)*1000"
$ws.Range("K127:K132").Formula = "=J127/(C127/1000)"

# Keep the sheet view pointed at the bottom of the now-longer data range
# (mirrors the freeze-pane scroll position Excel records after the user
# scrolls down to the newly entered rows).
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K135").Select()
